$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$rows = @(
    @{Row=1; Label="Danh mục lương"; Val=0},
    @{Row=2; Label="Tổng công tại CẦN THƠ"; Val=28},
    @{Row=3; Label="Lương cơ bản tại CẦN THƠ"; Val=7000000},
    @{Row=4; Label="Chiết khấu sale chính tại CẦN THƠ"; Val=0},
    @{Row=5; Label="Chiết khấu sale phụ tại CẦN THƠ"; Val=0},
    @{Row=6; Label="Đơn 1 bác sĩ tại CẦN THƠ"; Val=0},
    @{Row=7; Label="Đơn 2 bác sĩ tại CẦN THƠ"; Val=0},
    @{Row=8; Label="Công phụ phẫu 1 tại CẦN THƠ"; Val=0},
    @{Row=9; Label="Công phụ phẫu 2 tại CẦN THƠ"; Val=0},
    @{Row=10; Label="Ứng lương tại CẦN THƠ"; Val=-1000000},
    @{Row=11; Label="Tổng công tại LONG XUYÊN"; Val=0},
    @{Row=12; Label="Lương công tác tại LONG XUYÊN"; Val=0},
    @{Row=13; Label="Lương cơ bản tại LONG XUYÊN"; Val=$null},
    @{Row=14; Label="Chiết khấu sale chính tại LONG XUYÊN"; Val=0},
    @{Row=15; Label="Chiết khấu sale phụ tại LONG XUYÊN"; Val=0},
    @{Row=16; Label="Đơn 1 bác sĩ tại LONG XUYÊN"; Val=0},
    @{Row=17; Label="Đơn 2 bác sĩ tại LONG XUYÊN"; Val=0},
    @{Row=18; Label="Công phụ phẫu 1 tại LONG XUYÊN"; Val=0},
    @{Row=19; Label="Công phụ phẫu 2 tại LONG XUYÊN"; Val=0},
    @{Row=20; Label="Ứng lương tại LONG XUYÊN"; Val=-0},
    @{Row=21; Label="Tổng công tại SÓC TRĂNG"; Val=0},
    @{Row=22; Label="Lương công tác tại SÓC TRĂNG"; Val=0},
    @{Row=23; Label="Lương cơ bản tại SÓC TRĂNG"; Val=$null},
    @{Row=24; Label="Chiết khấu sale chính tại SÓC TRĂNG"; Val=0},
    @{Row=25; Label="Chiết khấu sale phụ tại SÓC TRĂNG"; Val=0},
    @{Row=26; Label="Đơn 1 bác sĩ tại SÓC TRĂNG"; Val=0},
    @{Row=27; Label="Đơn 2 bác sĩ tại SÓC TRĂNG"; Val=0},
    @{Row=28; Label="Công phụ phẫu 1 tại SÓC TRĂNG"; Val=0},
    @{Row=29; Label="Công phụ phẫu 2 tại SÓC TRĂNG"; Val=0},
    @{Row=30; Label="Ứng lương tại SÓC TRĂNG"; Val=-0},
    @{Row=31; Label="Tổng lương tại CẦN THƠ"; Val=6000000},
    @{Row=32; Label="Tổng lương tại LONG XUYÊN"; Val=0},
    @{Row=33; Label="Tổng lương tại SÓC TRĂNG"; Val=0},
    @{Row=34; Label="Tổng lương"; Val=6000000}
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Label
    $ws.Cells.Item($row.Row, 2).Value = $row.Val
}
